# Update the TestCases / TestData workbook: the old "LoginTest" +
# "ValidateCRMTest" (CRM) cases are replaced with two Google-search cases,
# "searchGoogle" and "searchGoogle2", each with its own block of test data.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("TestCases")
$ws2 = $wb.Worksheets.Item("TestData")

# ---- Sheet "TestCases" ------------------------------------------------
# A1/B1 header ("TestCases"/"Runmode") stay as-is.
$ws1.Range("A2").Value = "searchGoogle"
$ws1.Range("B2").Value = "Y"
$ws1.Range("A3").Value = "searchGoogle2"
$ws1.Range("B3").Value = "Y"

# ---- Sheet "TestData" --------------------------------------------------
# Block 1: searchGoogle
$ws2.Range("A1").Value = "searchGoogle"

$ws2.Range("A2").Value = "Runmode"
$ws2.Range("B2").Value = "SearchKeyword"
$ws2.Range("C2").Value = "password"
$ws2.Range("D2").Value = "browser"

$ws2.Range("A3").Value = "Y"
$ws2.Range("B3").Value = "Selenium Automation demo 1"
$ws2.Range("C3").Value = "Admin@123"
$ws2.Range("D3").Value = "chrome"

# Rows 4 and 5 remain blank (spacer between the two blocks).

# Block 2: searchGoogle2
$ws2.Range("A6").Value = "searchGoogle2"

$ws2.Range("A7").Value = "Runmode"
$ws2.Range("B7").Value = "SearchKeyword"
$ws2.Range("C7").Value = "password"
$ws2.Range("D7").Value = "browser"

$ws2.Range("A8").Value = "Y"
$ws2.Range("B8").Value = "Selenium Automation demo 2"
$ws2.Range("C8").Value = "Admin@123"
$ws2.Range("D8").Value = "chrome"

# ---- Selection / active-sheet bookkeeping -----------------------------
# TestCases becomes the active/selected tab (A4 selected); TestData keeps
# its own remembered selection (B2) but is no longer the active tab.
$ws1.Activate()
$ws1.Range("A4").Select()
$ws2.Range("B2").Select()
$ws1.Activate()
